# Weekly update: a new pricing report (date 2021-10-07 / serial 44476) for
# "Larga vida" Tomate at Vega Monumental Concepción is inserted as two new
# rows right before the existing row 231, pushing all the following rows
# (old 231..260) down by two (new 233..262).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above what is currently row 231 (shifts 231.. down to 233..)
$ws.Rows.Item(231).Insert()
$ws.Rows.Item(231).Insert()

# New row 231: Primera quality
$ws.Range("A231").Value = 11
$ws.Range("B231").Value = "Vega Monumental Concepción"
$ws.Range("C231").Value = "Bíobío"
$ws.Range("D231").Value = 44476
$ws.Range("E231").Value = 8
$ws.Range("F231").Value = 100112020
$ws.Range("G231").Value = "Tomate"
$ws.Range("H231").Value = "Larga vida"
$ws.Range("I231").Value = "Primera"
$ws.Range("J231").Value = 800
$ws.Range("K231").Value = 20000
$ws.Range("L231").Value = 21000
$ws.Range("M231").Value = 20500
$ws.Range("N231").Value = "`$/bandeja 18 kilos"
$ws.Range("O231").Value = "Región de Arica y Parinacota"
$ws.Range("P231").Value = 1139
$ws.Range("Q231").Value = 18
$ws.Range("R231").Value = "Hortaliza"

# New row 232: Segunda quality
$ws.Range("A232").Value = 11
$ws.Range("B232").Value = "Vega Monumental Concepción"
$ws.Range("C232").Value = "Bíobío"
$ws.Range("D232").Value = 44476
$ws.Range("E232").Value = 8
$ws.Range("F232").Value = 100112020
$ws.Range("G232").Value = "Tomate"
$ws.Range("H232").Value = "Larga vida"
$ws.Range("I232").Value = "Segunda"
$ws.Range("J232").Value = 400
$ws.Range("K232").Value = 19000
$ws.Range("L232").Value = 19000
$ws.Range("M232").Value = 19000
$ws.Range("N232").Value = "`$/bandeja 18 kilos"
$ws.Range("O232").Value = "Región de Arica y Parinacota"
$ws.Range("P232").Value = 1056
$ws.Range("Q232").Value = 18
$ws.Range("R232").Value = "Hortaliza"
